$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column header white color (matches the original "FFFFFFFF" run color used
# throughout row 1 of this report table).
$White = 16777215

function Set-HeaderBold($range) {
    # NOTE: calling Characters(1, <full length>) in a single shot is a
    # no-op in this COM host (it only sticks if the run actually gets
    # split), so we deliberately split the call in two pieces that
    # together cover the whole string; Excel then re-merges them into a
    # single run once both share identical formatting.
    $txt = $range.Value2
    $len = $txt.Length
    if ($len -gt 1) {
        $range.Characters(1, $len - 1).Font.Bold = $true
        $range.Characters($len, 1).Font.Bold = $true
    } else {
        $range.Characters(1, $len).Font.Bold = $true
    }
}

function Set-HeaderText($range, $text) {
    $range.Value = $text
    $len = $text.Length
    if ($len -gt 1) {
        $range.Characters(1, $len - 1).Font.Bold = $true
        $range.Characters($len, 1).Font.Bold = $true
        $range.Characters(1, $len - 1).Font.Color = $White
        $range.Characters($len, 1).Font.Color = $White
    } else {
        $range.Characters(1, $len).Font.Bold = $true
        $range.Characters(1, $len).Font.Color = $White
    }
}

# B1, D1 and F1 get new French wording (and become bold, like the rest of
# row 1's metric headers).
Set-HeaderText $ws.Range("B1") "Total des ventes de chaï (unités)"
Set-HeaderText $ws.Range("D1") "Ventes de chaï préconfectionné (unités)"
Set-HeaderText $ws.Range("F1") "Recherches en ligne de chaï"

# C1 and E1 keep their existing text, but also become bold.
Set-HeaderBold $ws.Range("C1")
Set-HeaderBold $ws.Range("E1")

# Keep the table's column headers (ListObject) in sync with the new cell
# text so the table definition and the worksheet stay consistent.
$table = $ws.ListObjects.Item(1)
$table.ListColumns.Item(2).Name = "Total des ventes de chaï (unités)"
$table.ListColumns.Item(4).Name = "Ventes de chaï préconfectionné (unités)"
$table.ListColumns.Item(6).Name = "Recherches en ligne de chaï"
